$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------------
# 1. Make room for the new daily-log rows: insert 7 blank rows before the
#    old row 19 ("Keyboard shortcuts" help block), which pushes that block
#    down to rows 26-28.
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Resize(7).Insert()

# ---------------------------------------------------------------------------
# 2. Column A dates.
#    Rows 12-18 become literal (non-formula) dates for the new daily logs.
#    Rows 4-11 keep the "previous + 7" formula, but the shared-formula group
#    is now restricted to A4:A11 (it used to stretch to A16).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 42686
$ws.Range("A13").Value = 42687
$ws.Range("A14").Value = 42688
$ws.Range("A15").Value = 42690
$ws.Range("A16").Value = 42691
$ws.Range("A17").Value = 42693
$ws.Range("A18").Value = 42696

$ws.Range("A4:A11").Formula = "=A3+7"

# ---------------------------------------------------------------------------
# 3. New daily-log content (rows 12-18).
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "• Worked on UI (SignIn, Teach, Settings)" + $nl + "• Sign in, sign up all works with databases , fixed navigation bars, got global variables to work"
$ws.Range("C12").Value = "• Worked on the database, made various functions. Also fixed errors"
$ws.Range("D12").Value = "• Worked on UI (Home, Learn) Created prototypes and statically populated names."
$ws.Range("E12").Value = "• Worked on the report" + $nl + "• Worked on UI communicated with the database" + $nl + "• Figured out how to code the transition of screens rather than using the Interface Builder"

$ws.Range("B13").Value = "• Worked on UI (SignIn, Teach, Settings)" + $nl + "• Worked on making the UI dyanmic" + $nl + "• Worked on populating courses in the search bar"
$ws.Range("C13").Value = "• Worked on the database, made various functions. Also fixed errors"
$ws.Range("D13").Value = "• Worked on UI (Home, Learn) Created prototypes and statically populated names." + $nl + "• Worked on making the UI dyanmic" + $nl + "• Worked on populating courses in the search bar"
$ws.Range("E13").Value = "• Worked on the report" + $nl + "• Worked on UI communicated with the database" + $nl + "• Figured out how to code the transition of screens rather than using the Interface Builder"

$ws.Range("B14").Value = "• Worked on UI (SignIn, Teach, Settings)" + $nl + "• Worked on making the UI dyanmic" + $nl + "• Worked on populating courses in the search bar"
$ws.Range("C14").Value = "• Worked on the database, made various functions. Also fixed errors"
$ws.Range("D14").Value = "• Worked on UI (Home, Learn)" + $nl + "• Worked on making the UI dyanmic" + $nl + "• Worked on populating courses in the search bar"
$ws.Range("E14").Value = "• Worked on UI communicated with the database" + $nl + "• Figured out how to code the transition of screens rather than using the Interface Builder"

# Row 15 stays blank (highlighted yellow placeholder).

$ws.Range("B16").Value = "• Finalized the App for presentation"
$ws.Range("C16").Value = "• Finalized the App for presentation and changed the UI "
$ws.Range("D16").Value = "• Finalized the App for presentation"
$ws.Range("E16").Value = "• Finalized the App for presentation"

$ws.Range("B17").Value = "• learnt how to colour pages properly and worked on prototyping for the HomeViewController"
$ws.Range("C17").Value = " -"
$ws.Range("D17").Value = " -"
$ws.Range("E17").Value = "• updated the SignInViewController to create Students as well as Users"

$ws.Range("B18").Value = "• Discusses the things remaining to work on"
$ws.Range("C18").Value = "• Discusses the things remaining to work on"
$ws.Range("D18").Value = "• Discusses the things remaining to work on"
$ws.Range("E18").Value = "• Discusses the things remaining to work on"

# Rows 19-20: blank placeholder rows (same look as the old empty weeks).
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""

# ---------------------------------------------------------------------------
# 4. Formatting.
# ---------------------------------------------------------------------------

# 4a. Column A: italic date format dd/mm/yyyy
$ws.Range("A1:A18").NumberFormat = "dd/mm/yyyy;@"
$ws.Range("A18").VerticalAlignment = -4108

# 4b. Column E width grows from 34.5 to 45, matching column D/B.
$ws.Columns.Item(5).ColumnWidth = 44.1

# 4c. Borders: every populated data cell in columns B-E gets a thin left
#     border, matching the rest of the sheet. Border edges are applied per
#     column (vertical ranges) so every row gets its own left edge.
$ws.Range("B12:B20").Borders.Item(7).LineStyle = 1
$ws.Range("C12:C20").Borders.Item(7).LineStyle = 1
$ws.Range("D12:D20").Borders.Item(7).LineStyle = 1
$ws.Range("E12:E20").Borders.Item(7).LineStyle = 1

# 4d. Wrap text + vertical centering for the new multi-line log rows.
$ws.Range("B12:D14").WrapText = $true
$ws.Range("B12:D14").VerticalAlignment = -4108
$ws.Range("E12:E14").WrapText = $true

$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("D16:E16").VerticalAlignment = -4108
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4108

$ws.Range("B17").WrapText = $true
$ws.Range("E17").WrapText = $true
$ws.Range("C17:D17").HorizontalAlignment = -4108

$ws.Range("B18:E18").WrapText = $true
$ws.Range("B18:E18").VerticalAlignment = -4108

# 4e. Row 15 highlighted yellow (placeholder week).
$ws.Range("B15:E15").Interior.Color = 65535

# 4f. Row heights for the new wrapped rows.
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 45

# ---------------------------------------------------------------------------
# 5. Selection matches the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("C18").Select()
